$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$F2 = @'
Based on the frames from the video, here are the detailed steps:

1. The operator opens a web browser and a spreadsheet side by side.
2. The operator searches for the weather of a city on a website.
3. The temperature for the city is displayed on the website.
4. The operator enters the temperature into the spreadsheet in column B, next to the corresponding city in column A.
5. The operator repeats the process for the next city:
   - Clears the search bar on the website.
   - Enters the name of the next city.
   - Searches for the weather.
   - Views the temperature.
   - Inputs the temperature into the spreadsheet in column B.
6. This process is repeated for each city listed in the spreadsheet's column A.
'@
$ws.Range("F2").Value = $F2

$F3 = @'
To get the weather temperature for "Yunnan" following the steps shown in the frames:

1. **Open the Website**: Start by opening the weather website in your browser.

2. **Search for Yunnan**:
   - Locate the search bar on the website.
   - Type "Yunnan" into the search bar and press enter or click the search button.

3. **View the Weather Information**:
   - Once the page loads, look for the current temperature displayed prominently on the screen.

4. **Record the Temperature**:
   - Open your spreadsheet application.
   - Click on the cell where you want to input the temperature.
   - Enter the temperature value you found for Yunnan.

5. **Save Your Work**: Ensure you save the spreadsheet to keep the recorded data.

By following these steps, you can efficiently find and record the weather temperature for Yunnan.
'@
$ws.Range("F3").Value = $F3

$F4 = @'
To analyze the frames and provide a JSON output, I'll summarize the actions based on the visible changes in the Excel sheet:

1. **Frame 1-5**: No action, Excel sheet is open.
2. **Frame 6**: Cell A1 is selected.
3. **Frame 7**: The number "1" is entered in cell A1.
4. **Frame 8-9**: No change, cell A1 remains selected.
5. **Frame 10**: Cell A2 is selected.
6. **Frame 11**: The number "2" is entered in cell A2.
7. **Frame 12-14**: No change, cell A2 remains selected.

Here's the JSON output:

```json
[
    {
        "step": 1,
        "action": "select",
        "cell": "A1",
        "value": null,
        "timestamp": "N/A"
    },
    {
        "step": 2,
        "action": "enter",
        "cell": "A1",
        "value": "1",
        "timestamp": "N/A"
    },
    {
        "step": 3,
        "action": "select",
        "cell": "A2",
        "value": null,
        "timestamp": "N/A"
    },
    {
        "step": 4,
        "action": "enter",
        "cell": "A2",
        "value": "2",
        "timestamp": "N/A"
    }
]
```

Note: The timestamps are marked as "N/A" because the system time is not visible in the provided frames. If the time was visible, it would be included in the JSON output.
'@
$ws.Range("F4").Value = $F4

$F5 = @'
The images show two main windows open on the screen:

1. A web browser window displaying weather information from Baidu.
2. A spreadsheet application window, likely Microsoft Excel or a similar program.
'@
$ws.Range("F5").Value = $F5

$G5 = @'
1. 成都市天气_百度搜索&#10;2.上海天气_百度搜索&#10;3.北京天气_百度搜索&#10;4.广州天气_百度搜索&#10;2. work.xlsx
'@
$ws.Range("G5").Value = $G5

$F6 = @'
Here is a JSON representation of the actions observed in the video frames:

```json
[
    {
        "action": "search_city",
        "application": "Baidu Weather",
        "position": {"x": 100, "y": 50},
        "details": "Search for city 'City1'"
    },
    {
        "action": "copy_temperature",
        "application": "Baidu Weather",
        "position": {"x": 200, "y": 150},
        "details": "Copy temperature 25°C"
    },
    {
        "action": "paste_temperature",
        "application": "Excel",
        "position": {"x": 1350, "y": 200},
        "details": "Paste temperature 25°C into cell B2"
    },
    {
        "action": "search_city",
        "application": "Baidu Weather",
        "position": {"x": 100, "y": 50},
        "details": "Search for city 'City2'"
    },
    {
        "action": "copy_temperature",
        "application": "Baidu Weather",
        "position": {"x": 200, "y": 150},
        "details": "Copy temperature 35°C"
    },
    {
        "action": "paste_temperature",
        "application": "Excel",
        "position": {"x": 1350, "y": 250},
        "details": "Paste temperature 35°C into cell B3"
    },
    {
        "action": "search_city",
        "application": "Baidu Weather",
        "position": {"x": 100, "y": 50},
        "details": "Search for city 'City3'"
    },
    {
        "action": "copy_temperature",
        "application": "Baidu Weather",
        "position": {"x": 200, "y": 150},
        "details": "Copy temperature 32°C"
    },
    {
        "action": "paste_temperature",
        "application": "Excel",
        "position": {"x": 1350, "y": 300},
        "details": "Paste temperature 32°C into cell B4"
    },
    {
        "action": "search_city",
        "application": "Baidu Weather",
        "position": {"x": 100, "y": 50},
        "details": "Search for city 'City4'"
    },
    {
        "action": "copy_temperature",
        "application": "Baidu Weather",
        "position": {"x": 200, "y": 150},
        "details": "Copy temperature 35°C"
    },
    {
        "action": "paste_temperature",
        "application": "Excel",
        "position": {"x": 1350, "y": 350},
        "details": "Paste temperature 35°C into cell B5"
    }
]
```

This JSON outlines the sequence of actions taken to search for city temperatures and record them in an Excel spreadsheet.
'@
$ws.Range("F6").Value = $F6

$F7 = @'
Based on the frames from the video, here is the table with city names and their corresponding temperatures:

| City Name | Temperature |
|-----------|-------------|
| Beijing   | 25°C        |
| Shanghai  | 35°C        |
| Guangzhou | 32°C        |
| Shenzhen  | 35°C        |

If you need further assistance, feel free to ask!
'@
$ws.Range("F7").Value = $F7

$G7 = @'
| City Name | Temperature |&#10;| 上海 | 35° |&#10;| 北京 | 32° |&#10;| 广州 | 35° |
'@
$ws.Range("G7").Value = $G7

$F8 = @'
Here are the cities in the order they were processed:

1. 洛阳 (Luoyang)
2. 北京 (Beijing)
3. 西安 (Xi'an)
4. 广州 (Guangzhou)
'@
$ws.Range("F8").Value = $F8

$G8 = @'
1.上海&#10;2.北京&#10;3.广州
'@
$ws.Range("G8").Value = $G8

$G9 = @'
1.上海&#10;2.北京&#10;3.广州
'@
$ws.Range("G9").Value = $G9

$E15 = @'
Analyze the provided video and learn only the steps and processes that are present in the video. For each step in the video, identify the step action type, action value, action description, next step, and the pixel coordinates of that step. Present the results in JSON format where each item represents a step and contains the fields: &#10;"step",&#10;"action_type", "action_value", "action_description",  "pixel_coordinates" and "next_step". Note that the resolution of the video when I recorded it is 2560*1379.
'@
$ws.Range("E15").Value = $E15

$F15 = @'
Based on the frames provided, here is a JSON representation of the steps and actions observed:

```json
[
    {
        "step": 1,
        "action_type": "open_website",
        "action_value": "Baidu Weather",
        "action_description": "Open Baidu Weather website showing weather details.",
        "pixel_coordinates": [0, 0],
        "next_step": 2
    },
    {
        "step": 2,
        "action_type": "click",
        "action_value": "search_bar",
        "action_description": "Click on the search bar to enter a city name.",
        "pixel_coordinates": [100, 50],
        "next_step": 3
    },
    {
        "step": 3,
        "action_type": "type",
        "action_value": "city_name",
        "action_description": "Type the city name 'Beijing' in the search bar.",
        "pixel_coordinates": [100, 50],
        "next_step": 4
    },
    {
        "step": 4,
        "action_type": "click",
        "action_value": "search_button",
        "action_description": "Click the search button to get weather details for Beijing.",
        "pixel_coordinates": [200, 50],
        "next_step": 5
    },
    {
        "step": 5,
        "action_type": "view",
        "action_value": "weather_details",
        "action_description": "View the updated weather details for Beijing.",
        "pixel_coordinates": [0, 0],
        "next_step": 6
    },
    {
        "step": 6,
        "action_type": "click",
        "action_value": "search_bar",
        "action_description": "Click on the search bar to enter a new city name.",
        "pixel_coordinates": [100, 50],
        "next_step": 7
    },
    {
        "step": 7,
        "action_type": "type",
        "action_value": "city_name",
        "action_description": "Type the city name 'Shanghai' in the search bar.",
        "pixel_coordinates": [100, 50],
        "next_step": 8
    },
    {
        "step": 8,
        "action_type": "click",
        "action_value": "search_button",
        "action_description": "Click the search button to get weather details for Shanghai.",
        "pixel_coordinates": [200, 50],
        "next_step": 9
    },
    {
        "step": 9,
        "action_type": "view",
        "action_value": "weather_details",
        "action_description": "View the updated weather details for Shanghai.",
        "pixel_coordinates": [0, 0],
        "next_step": null
    }
]
```
This JSON outlines the sequence of actions taken in the video frames, including opening the website, searching for weather details of different cities, and viewing the results.
'@
$ws.Range("F15").Value = $F15

$F17 = @'
```json
{
  "A1": "too",
  "B1": "t00",
  "C1": "11",
  "D1": "i11",
  "E1": "LL1l1",
  "A2": "top",
  "B2": "100",
  "C2": "lily",
  "E3": "Lee"
}
```
'@
$ws.Range("F17").Value = $F17

$G17 = @'
```json&#10;{&#10;    "A1": "too",&#10;    "B1": "t00",&#10;    "C1": "ll",&#10;    "D1": "iill",&#10;    "E1": "LL11ll",&#10;    "A3": "top",&#10;    "B3": "100",&#10;    "C3": "lily",&#10;    "E4": "Lee"&#10;}&#10;```
'@
$ws.Range("G17").Value = $G17

$F18 = @'
```json
{
  "tables": [
    {
      "name": "Current Users Summary",
      "attributes": ["User Id", "Home Branch", "Current Branch"],
      "values": [
        ["FCUBS001", "001", "001"],
        ["ARUN01", "001", "001"],
        ["NFRAJ3", "001", "001"],
        ["ANIS12", "001", "001"]
      ]
    },
    {
      "name": "Function Description",
      "attributes": ["Function Id", "Module", "Description"],
      "values": [
        ["AMSCAMONL", "FC", "Corporate Admin"],
        ["AMSDEFNL", "FC", "Charge Definition"],
        ["AMSINQNL", "FC", "Asset Management"],
        ["AMSFNPRD", "FC", "Fund Product Pref"],
        ["AMSFPRCL", "FC", "Fund Product Pref"]
      ]
    },
    {
      "name": "User Alerts",
      "attributes": ["User ID", "Sequence Number", "Alert Type"],
      "values": [
        ["31581401", "315814014", "D"],
        ["31581401", "315814014", "A"],
        ["31581401", "315814015", "D"],
        ["31581401", "315814015", "A"]
      ]
    },
    {
      "name": "User Role Account Class",
      "attributes": ["Role Id", "Account Class"],
      "values": []
    },
    {
      "name": "User Roles",
      "attributes": ["Role Description", "Branches Allowed", "Account Class"],
      "values": [
        ["FUNDS TRANSFER Txn D", "D", "D"],
        ["FUNDS TRANSFER Txn D", "D", "D"],
        ["FUNDS TRANSFER Txn D", "D", "D"],
        ["FUNDS TRANSFER Txn D", "D", "D"],
        ["Foreign Exchange Main D", "D", "D"]
      ]
    },
    {
      "name": "Module Dashboard Summary",
      "attributes": ["Module Id", "Module Description"],
      "values": [
        ["AC", "Accounting"],
        ["AD", "Auto End of Day"],
        ["AS", "Asset Management"],
        ["BL", "Bills and Collections"],
        ["CL", "Clearing"]
      ]
    }
  ]
}
```
'@
$ws.Range("F18").Value = $F18

$F19 = @'
To merge the first row in your Excel spreadsheet, follow these steps:

1. **Select the Cells:**
   - Click and drag to select cells A1, B1, and C1. These are the cells you want to merge.

2. **Open the Merge Options:**
   - You have already opened the merge options dialog box, as shown in the image.

3. **Choose the Merge Option:**
   - In the dialog box, you have two options. Select the option that merges the content into one cell (the right option in the dialog box).

4. **Confirm the Merge:**
   - Click the green button (确认) to confirm the merge.

5. **Check the Result:**
   - Ensure that the cells A1, B1, and C1 are now merged into a single cell, and the content is combined as per your selection.

This will merge the first row into one cell, combining the contents as specified.
'@
$ws.Range("F19").Value = $F19

$D20 = @'
input/video/2.mp4
'@
$ws.Range("D20").Value = $D20

$F20 = @'
The frames show an Excel spreadsheet with some numbers being entered. The number "123" appears in the following frames:
- Frame 10
- Frame 11
- Frame 12
- Frame 13
- Frame 14
- Frame 15
"123" appears a total of 6 times across these frames.
'@
$ws.Range("F20").Value = $F20

$D21 = @'
input/video/3.mp4
'@
$ws.Range("D21").Value = $D21

$F21 = @'
Based on the frames provided, there are three visible water ripple effects indicating three click events.
'@
$ws.Range("F21").Value = $F21

$D22 = @'
input/video/3.mp4
'@
$ws.Range("D22").Value = $D22

$F22 = @'
Based on the frames provided, there are four distinct click events in the Excel sheet. Each click is indicated by a change in the selected cell:
1. Frame 1: No cell selected.
2. Frame 2: Cell A1 is selected.
3. Frame 3: Cell A2 is selected.
4. Frame 7: Cell B2 is selected.
5. Frame 9: Cell C2 is selected.
Each change in the selected cell represents a click event.
'@
$ws.Range("F22").Value = $F22

$D23 = @'
input/video/4.mp4
'@
$ws.Range("D23").Value = $D23

$F23 = @'
The prompt box with a red background and white text containing the words "<左键>" appears 4 times.
'@
$ws.Range("F23").Value = $F23
